$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (Model, Configuration, F1, Accuracy, Precision, Recall)
$cfg1 = "CV + tfidf + ngram(1)"
$cfg2 = "CV + tfidf + ngram(2)"
$cfg3 = "CV + tfidf + ngram(3)"

$data = @(
    @("Logistic Regression",      $cfg1, 88.91, 82.66, 81.66, 97.80),
    @("Multinomial Naive Bayes",  $cfg1, 88.44, 81.48, 79.57, 99.54),
    @("Support Vector Machines",  $cfg1, 83.94, 72.78, 72.78, 100),
    @("Decision Tree",            $cfg1, 88.55, 82.58, 83.43, 94.59),
    @("Random Forest",            $cfg1, 88.37, 82.25, 83.41, 94.22),

    @("Logistic Regression",      $cfg2, 88.75, 82.13, 80.63, 98.96),
    @("Multinomial Naive Bayes",  $cfg2, 88.19, 80.99, 79.25, 99.72),
    @("Support Vector Machines",  $cfg2, 83.94, 72.78, 72.78, 100),
    @("Decision Tree",            $cfg2, 87.78, 81.39, 82.87, 93.60),
    @("Random Forest",            $cfg2, 87.76, 81.35, 82.61, 93.78),

    @("Logistic Regression",      $cfg3, 88.44, 81.56, 80.06, 99.06),
    @("Multinomial Naive Bayes",  $cfg3, 88.07, 80.74, 78.98, 99.82),
    @("Support Vector Machines",  $cfg3, 83.94, 72.78, 72.78, 100),
    @("Decision Tree",            $cfg3, 87.36, 80.66, 82.34, 93.30),
    @("Random Forest",            $cfg3, 87.54, 81.02, 82.48, 93.50)
)

$row = 12
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}

# Mirror the final view state: scrolled down with C27 selected/active.
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 7 } catch {}
$ws.Range("C27").Select() | Out-Null
